# Issue 84: Create content outline for feature documentation pages Update Issue 84
#
# Changes applied:
#  1. Slide master + all 11 slide layouts: the cached text of the
#     "Date Placeholder" auto-update date field changes from 12/14/2013
#     to 12/24/2013.
#  2. Slide master shape "Rectangle 6": its line/border color changes
#     from a fixed RGB (0070C0) to the theme color White/Background 1,
#     Darker 15% (schemeClr bg1 + lumMod 85000 == sRGB D9D9D9).
#  3. Slide 1 ("Create better slides with less effort" title slide) is
#     hidden from the slide show.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- 1. Update the cached date field text (master + every layout) ---
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "12/24/2013"

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Set-DatePlaceholderText $layout.Shapes "12/24/2013"
}

# --- 2. Recolor the master's border rectangle line to bg1 / darker 15% ---
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 6") {
        $shp.Line.ForeColor.RGB = 14277081   # 0xD9D9D9 == bg1 lumMod 85000
    }
}

# --- 3. Hide slide 1 from the slide show ---
$s1 = $p.Slides.Item(1)
$s1.SlideShowTransition.Hidden = 1
